$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = 4
    3 = 8
    4 = 2
    5 = 1
    6 = 6
    7 = 4
    8 = 3
    9 = 8
    10 = 1
    11 = 7
    12 = 7
    13 = 8
    14 = 8
    15 = 6
    16 = 5
    17 = 5
    18 = 7
    19 = 9
    20 = 6
    21 = 8
    22 = 5
    23 = 2
    24 = 7
    25 = 4
    26 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
